# Rewrite the "KEY ACHIEVEMENTS AND IMPACT" bullet list to impact-focused
# accomplishment statements, and trim it from six bullets down to four.
#
# The same underlying sentences (e.g. "Achieved 87% prediction accuracy...")
# also appear verbatim in the "PROFESSIONAL EXPERIENCE" section earlier in
# the document, so we must not rely on a document-wide Find/Replace. Instead
# we anchor on the "KEY ACHIEVEMENTS AND IMPACT" heading and walk forward
# from there to find the exact paragraphs to edit/remove.

$d = $word.ActiveDocument

# Locate the "Impact" sub-heading paragraph that immediately follows the
# "KEY ACHIEVEMENTS AND IMPACT" section heading.
$anchor = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text.Trim()
    if ($t -eq "KEY ACHIEVEMENTS AND IMPACT") {
        $anchor = $i
        break
    }
}

if ($anchor -eq $null) {
    throw "Could not find 'KEY ACHIEVEMENTS AND IMPACT' heading"
}

# The six bullet paragraphs sit right after the heading + "Impact" sub-head.
$b1 = $anchor + 2
$b2 = $anchor + 3
$b3 = $anchor + 4
$b4 = $anchor + 5
$b5 = $anchor + 6
$b6 = $anchor + 7

# Sanity-check the original text before mutating anything.
$expect1 = "• Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%"
$expect2 = "• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from " + [char]0x00B1 + "4.2% to " + [char]0x00B1 + "2.1%"
$expect3 = "• Built cloud-based data warehouse solutions on AWS processing billions of records with 99.94% accuracy"
$expect4 = "• Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations"
$expect5 = "• Developed longitudinal data analysis methods using geospatial techniques that improved segmentation accuracy by 34% and survey incidence rates by 28%, reducing polling costs while increasing response quality"
$expect6 = "• Designed ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial datasets"

if ($d.Paragraphs.Item($b1).Range.Text.Trim() -ne $expect1) { throw "bullet1 mismatch: " + $d.Paragraphs.Item($b1).Range.Text }
if ($d.Paragraphs.Item($b2).Range.Text.Trim() -ne $expect2) { throw "bullet2 mismatch: " + $d.Paragraphs.Item($b2).Range.Text }
if ($d.Paragraphs.Item($b3).Range.Text.Trim() -ne $expect3) { throw "bullet3 mismatch: " + $d.Paragraphs.Item($b3).Range.Text }
if ($d.Paragraphs.Item($b4).Range.Text.Trim() -ne $expect4) { throw "bullet4 mismatch: " + $d.Paragraphs.Item($b4).Range.Text }
if ($d.Paragraphs.Item($b5).Range.Text.Trim() -ne $expect5) { throw "bullet5 mismatch: " + $d.Paragraphs.Item($b5).Range.Text }
if ($d.Paragraphs.Item($b6).Range.Text.Trim() -ne $expect6) { throw "bullet6 mismatch: " + $d.Paragraphs.Item($b6).Range.Text }

# Rewrite the three kept bullets in place (preserves paragraph formatting).
$d.Paragraphs.Item($b1).Range.Text = "• Predictive excellence: Achieved 87% voter turnout accuracy vs. 71% industry standard"
$d.Paragraphs.Item($b2).Range.Text = "• Reduced polling margins from " + [char]0x00B1 + "4.2% to " + [char]0x00B1 + "2.1%"
$d.Paragraphs.Item($b3).Range.Text = "• Methodological advancement: Improved segmentation accuracy 34% and survey incidence 28%"
$d.Paragraphs.Item($b5).Range.Text = "• Reduced polling costs while increasing quality"

# Remove the two dropped bullets entirely (including their paragraph marks).
# Delete from the higher index down so the lower index stays valid.
$d.Paragraphs.Item($b6).Range.Delete()
$d.Paragraphs.Item($b4).Range.Delete()
